$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date-formatted style (from H2:I2) to reuse for new date cells
$dateStyleSource = $ws.Range("H2:I2")

# Row 55
$ws.Range("A55").Value = 141
$ws.Range("B55").Value = "Linked List Cycle"
$ws.Range("C55").Value = "#linked-list #two-pointers #必背 "
$ws.Range("D55").Value = "easy"
$ws.Range("E55").Value = 2
$ws.Range("F55").Value = 3
$ws.Range("G55").Value = 20
$ws.Range("H55").Value = 45842
$ws.Range("I55").Value = 45842
$ws.Range("J55").Value = "?"
$dateStyleSource.Copy()
$ws.Range("H55:I55").PasteSpecial(-4122)
$ws.Rows.Item(55).RowHeight = 34

# Row 56
$ws.Range("A56").Value = 160
$ws.Range("B56").Value = "Intersection of Two Linked Lists"
$ws.Range("C56").Value = "#linked-list #two-pointers #核心 "
$ws.Range("D56").Value = "easy"
$ws.Range("E56").Value = 3
$ws.Range("F56").Value = 2
$ws.Range("G56").Value = 10
$ws.Range("H56").Value = 45842
$ws.Range("I56").Value = 45842
$dateStyleSource.Copy()
$ws.Range("H56:I56").PasteSpecial(-4122)
$ws.Rows.Item(56).RowHeight = 34

# Row 57
$ws.Range("A57").Value = 234
$ws.Range("B57").Value = "Palindrome Linked List"
$ws.Range("C57").Value = "#linked-list #two-pointers #核心 "
$ws.Range("D57").Value = "easy"
$ws.Range("E57").Value = 3
$ws.Range("F57").Value = 2
$ws.Range("G57").Value = 20
$ws.Range("H57").Value = 45842
$ws.Range("I57").Value = 45842
$ws.Range("J57").Value = "?"
$dateStyleSource.Copy()
$ws.Range("H57:I57").PasteSpecial(-4122)
$ws.Rows.Item(57).RowHeight = 34

# Row 58
$ws.Range("A58").Value = 1394
$ws.Range("B58").Value = "Find Lucky Integer in an Array"
$ws.Range("D58").Value = "easy"
$ws.Range("E58").Value = 1
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 45843
$ws.Range("I58").Value = 45843
$dateStyleSource.Copy()
$ws.Range("H58:I58").PasteSpecial(-4122)
$ws.Rows.Item(58).RowHeight = 34

# Row 59
$ws.Range("A59").Value = 1865
$ws.Range("B59").Value = "Finding Pairs With a Certain Sum"
$ws.Range("C59").Value = "#array #hash-table #design"
$ws.Range("D59").Value = "medium"
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 24
$ws.Range("H59").Value = 45844
$ws.Range("I59").Value = 45844
$dateStyleSource.Copy()
$ws.Range("H59:I59").PasteSpecial(-4122)
$ws.Rows.Item(59).RowHeight = 34

# Row 60
$ws.Range("A60").Value = 1353
$ws.Range("B60").Value = "Maximum Number of Events That Can Be Attended"
$ws.Range("C60").Value = "#array #greedy #sorting #heap "
$ws.Range("D60").Value = "medium"
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 28
$ws.Range("H60").Value = 45845
$ws.Range("I60").Value = 45845
$dateStyleSource.Copy()
$ws.Range("H60:I60").PasteSpecial(-4122)
$ws.Rows.Item(60).RowHeight = 51

# Row 61
$ws.Range("A61").Value = 328
$ws.Range("B61").Value = "Odd Even Linked List"
$ws.Range("C61").Value = "#linked-list"
$ws.Range("D61").Value = "medium"
$ws.Range("E61").Value = 1
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 20
$ws.Range("H61").Value = 45845
$ws.Range("I61").Value = 45845
$ws.Range("J61").Value = "?"
$dateStyleSource.Copy()
$ws.Range("H61:I61").PasteSpecial(-4122)
$ws.Rows.Item(61).RowHeight = 17

$ws.Range("J61").Select()

$excel.Application.CutCopyMode = $false
